$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing existing rows 5-14 down to 6-15.
$ws.Rows(5).Insert()

# Copy style (date format) from the cell above into the new D5 cell.
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row 5 values (same as the other rows, new date + prices).
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 44525
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 300000000
$ws.Range("G5").Value = "Espárragos"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 360
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 900
$ws.Range("M5").Value = 850
$ws.Range("N5").Value = "$/kilo"
$ws.Range("O5").Value = "Provincia de Diguillín"
$ws.Range("P5").Value = 850
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
